$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 3 describes the Marine benthic input dataset. The workflow now
# points at the processed "Marine_Benthic" intermediate raster rather
# than the raw "Natural Values Ecosystems" input, so update the raw
# data path and dataset name accordingly.
$ws.Range("B3").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\Marine_Benthic.tif"
$ws.Range("A3").Value = "Marine"

# Move the active selection to A4, matching where the author left off.
$ws.Range("A4").Select()
